$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "G1"
$ws.Range("B2").Value = "Test1"
$ws.Range("C2").Value = 45889
$ws.Range("C2").NumberFormat = "yyyy-mm-dd"
$ws.Range("C2").NumberFormat = "YYYY-MM-DD"
$ws.Range("D2").Value = 1.01
$ws.Range("E2").Value = 100
$ws.Range("F2").Value = 0.01
